# Apply "Bug fixes in albedo training" edit to "Maps 2 RGB" sheet:
# - Restyle existing rows 25-26 (B:I) to the new red-font 0.0000 style
# - Fill in missing data for rows 27-28 (V4.06.7 / V4.06.8) using the same new style
# - Add a brand-new data row 29 (V4.07.5) with the new style
# - Add new label-only rows 30-32 (V4.07.6, V4.07.7, V4.07.8)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Maps 2 RGB")

# New style: number format 0.0000, red font (matches new cellXf / font added to styles.xml)
$dataRange = $ws.Range("B25:I29")
$dataRange.NumberFormat = "0.0000"
$dataRange.Font.Color = 255

# Row 25 (V4.06.5) - values unchanged, only the style changes
$ws.Cells.Item(25,2).Value = 17.706900000000001
$ws.Cells.Item(25,3).Value = 0.76649999999999996
$ws.Cells.Item(25,4).Value = 19.653099999999998
$ws.Cells.Item(25,5).Value = 0.86029999999999995
$ws.Cells.Item(25,6).Value = 22.648900000000001
$ws.Cells.Item(25,7).Value = 0.90410000000000001
$ws.Cells.Item(25,8).Value = 21.154800000000002
$ws.Cells.Item(25,9).Value = 0.85070000000000001

# Row 26 (V4.06.6) - values unchanged, only the style changes
$ws.Cells.Item(26,2).Value = 17.365400000000001
$ws.Cells.Item(26,3).Value = 0.76890000000000003
$ws.Cells.Item(26,4).Value = 19.2791
$ws.Cells.Item(26,5).Value = 0.85189999999999999
$ws.Cells.Item(26,6).Value = 23.146000000000001
$ws.Cells.Item(26,7).Value = 0.90639999999999998
$ws.Cells.Item(26,8).Value = 22.748899999999999
$ws.Cells.Item(26,9).Value = 0.86929999999999996

# Row 27 (V4.06.7) - new values
$ws.Cells.Item(27,2).Value = 12.3668
$ws.Cells.Item(27,3).Value = 0.71499999999999997
$ws.Cells.Item(27,4).Value = 18.418500000000002
$ws.Cells.Item(27,5).Value = 0.84709999999999996
$ws.Cells.Item(27,6).Value = 22.5608
$ws.Cells.Item(27,7).Value = 0.86829999999999996
$ws.Cells.Item(27,8).Value = 16.315799999999999
$ws.Cells.Item(27,9).Value = 0.82110000000000005

# Row 28 (V4.06.8) - new values
$ws.Cells.Item(28,2).Value = 15.368399999999999
$ws.Cells.Item(28,3).Value = 0.74960000000000004
$ws.Cells.Item(28,4).Value = 19.710999999999999
$ws.Cells.Item(28,5).Value = 0.877
$ws.Cells.Item(28,6).Value = 23.2363
$ws.Cells.Item(28,7).Value = 0.88929999999999998
$ws.Cells.Item(28,8).Value = 20.630400000000002
$ws.Cells.Item(28,9).Value = 0.87549999999999994

# Row 29 (V4.07.5) - brand new row
$ws.Cells.Item(29,1).Value = "V4.07.5"
$ws.Cells.Item(29,2).Value = 18.610700000000001
$ws.Cells.Item(29,3).Value = 0.79330000000000001
$ws.Cells.Item(29,4).Value = 20.399799999999999
$ws.Cells.Item(29,5).Value = 0.88759999999999994
$ws.Cells.Item(29,6).Value = 20.9389
$ws.Cells.Item(29,7).Value = 0.86960000000000004
$ws.Cells.Item(29,8).Value = 25.526199999999999
$ws.Cells.Item(29,9).Value = 0.91520000000000001

# Rows 30-32: new label-only rows
$ws.Cells.Item(30,1).Value = "V4.07.6"
$ws.Cells.Item(31,1).Value = "V4.07.7"
$ws.Cells.Item(32,1).Value = "V4.07.8"

# Update the sheet view to match (scrolled down a bit further, selection moved to J27)
$ws.Application.ActiveWindow.ScrollRow = 16
[void]$ws.Range("J27").Select()
